$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 15
$ws_ALC.Range("H15").Value = 1056.4193
$ws_ALC.Range("I15").Value = 1056.4193
$ws_ALC.Range("K15").Value = 3169.2579
$ws_ALC.Range("M15").Value = -3000.2579

# ALC row 33
$ws_ALC.Range("H33").Value = 105.4
$ws_ALC.Range("J33").Value = 125
$ws_ALC.Range("L33").Value = 125
$ws_ALC.Range("N33").Value = -583

# ALC row 39
$ws_ALC.Range("H39").Value = 312.72726
$ws_ALC.Range("I39").Value = 55.857143
$ws_ALC.Range("J39").Value = 762.25
$ws_ALC.Range("K39").Value = 167.571429
$ws_ALC.Range("L39").Value = 2286.75
$ws_ALC.Range("M39").Value = 128.428571
$ws_ALC.Range("N39").Value = -2878.75

# ALC row 43
$ws_ALC.Range("H43").Value = 1279.6471
$ws_ALC.Range("I43").Value = 943
$ws_ALC.Range("J43").Value = 1383.2307
$ws_ALC.Range("K43").Value = 943
$ws_ALC.Range("L43").Value = 1383.2307
$ws_ALC.Range("M43").Value = -874
$ws_ALC.Range("N43").Value = -1521.2307

# ALC row 51
$ws_ALC.Range("H51").Value = 5070.3335
$ws_ALC.Range("J51").Value = 5188.4
$ws_ALC.Range("L51").Value = 5188.4
$ws_ALC.Range("N51").Value = -6156.4

# ALC row 55
$ws_ALC.Range("H55").Value = 330.1111
$ws_ALC.Range("I55").Value = 276
$ws_ALC.Range("J55").Value = 397.75
$ws_ALC.Range("K55").Value = 276
$ws_ALC.Range("L55").Value = 397.75
$ws_ALC.Range("M55").Value = -62
$ws_ALC.Range("N55").Value = -825.75

# ALC row 58
$ws_ALC.Range("H58").Value = 1126.5834
$ws_ALC.Range("J58").Value = 3460
$ws_ALC.Range("L58").Value = 10380
$ws_ALC.Range("N58").Value = -10680

# ALC row 62
$ws_ALC.Range("H62").Value = 2474.5
$ws_ALC.Range("J62").Value = 4000
$ws_ALC.Range("L62").Value = 4000
$ws_ALC.Range("N62").Value = -5248

# ALC row 65
$ws_ALC.Range("H65").Value = 2474.5
$ws_ALC.Range("J65").Value = 4000
$ws_ALC.Range("L65").Value = 20000
$ws_ALC.Range("N65").Value = -26240

# ALC row 100
$ws_ALC.Range("H100").Value = 5500
$ws_ALC.Range("J100").Value = 10000
$ws_ALC.Range("L100").Value = 10000
$ws_ALC.Range("N100").Value = -11082

# ALC row 138
$ws_ALC.Range("H138").Value = 1686.25
$ws_ALC.Range("I138").Value = 1608.4
$ws_ALC.Range("K138").Value = 4825.200000000001
$ws_ALC.Range("M138").Value = 314.7999999999993

# ARM row 32
$ws_ARM.Range("H32").Value = 5272.5386
$ws_ARM.Range("I32").Value = 3715.0715
$ws_ARM.Range("K32").Value = 3715.0715
$ws_ARM.Range("M32").Value = -3428.0715

# ARM row 92
$ws_ARM.Range("H92").Value = 50000
$ws_ARM.Range("J92").Value = 50000
$ws_ARM.Range("L92").Value = 50000
$ws_ARM.Range("N92").Value = -54992

# ARM row 140
$ws_ARM.Range("H140").Value = 0
$ws_ARM.Range("J140").Value = 0
$ws_ARM.Range("L140").Value = 0
$ws_ARM.Range("N140").ClearContents()

# BSM row 24
$ws_BSM.Range("H24").Value = 0
$ws_BSM.Range("I24").Value = 0
$ws_BSM.Range("K24").Value = 0
$ws_BSM.Range("M24").ClearContents()

# BSM row 25
$ws_BSM.Range("H25").Value = 25338.666
$ws_BSM.Range("I25").Value = 0
$ws_BSM.Range("J25").Value = 25338.666
$ws_BSM.Range("K25").Value = 0
$ws_BSM.Range("L25").Value = 25338.666
$ws_BSM.Range("M25").ClearContents()
$ws_BSM.Range("N25").Value = -25808.666

# BSM row 76
$ws_BSM.Range("H76").Value = 29666
$ws_BSM.Range("J76").Value = 29666
$ws_BSM.Range("L76").Value = 29666
$ws_BSM.Range("N76").Value = -30296

# BSM row 79
$ws_BSM.Range("H79").Value = 29666
$ws_BSM.Range("J79").Value = 29666
$ws_BSM.Range("L79").Value = 29666
$ws_BSM.Range("N79").Value = -31850

# BSM row 95
$ws_BSM.Range("H95").Value = 71896
$ws_BSM.Range("J95").Value = 71896
$ws_BSM.Range("L95").Value = 71896
$ws_BSM.Range("N95").Value = -77388

# BSM row 122
$ws_BSM.Range("H122").Value = 68000
$ws_BSM.Range("J122").Value = 68000
$ws_BSM.Range("L122").Value = 68000
$ws_BSM.Range("N122").Value = -77800

# CRP row 16
$ws_CRP.Range("H16").Value = 971
$ws_CRP.Range("I16").Value = 850
$ws_CRP.Range("J16").Value = 1213
$ws_CRP.Range("K16").Value = 850
$ws_CRP.Range("L16").Value = 1213
$ws_CRP.Range("M16").Value = -563
$ws_CRP.Range("N16").Value = -1787

# CRP row 62
$ws_CRP.Range("H62").Value = 0
$ws_CRP.Range("I62").Value = 0
$ws_CRP.Range("K62").Value = 0
$ws_CRP.Range("M62").ClearContents()

# CRP row 65
$ws_CRP.Range("H65").Value = 0
$ws_CRP.Range("I65").Value = 0
$ws_CRP.Range("K65").Value = 0
$ws_CRP.Range("M65").ClearContents()

# CRP row 103
$ws_CRP.Range("H103").Value = 0
$ws_CRP.Range("I103").Value = 0
$ws_CRP.Range("K103").Value = 0
$ws_CRP.Range("M103").ClearContents()

# CRP row 113
$ws_CRP.Range("H113").Value = 971
$ws_CRP.Range("I113").Value = 850
$ws_CRP.Range("J113").Value = 1213
$ws_CRP.Range("K113").Value = 850
$ws_CRP.Range("L113").Value = 1213
$ws_CRP.Range("M113").Value = 1320
$ws_CRP.Range("N113").Value = -5553

# CRP row 122
$ws_CRP.Range("H122").Value = 1517.3462
$ws_CRP.Range("I122").Value = 1675.5714
$ws_CRP.Range("J122").Value = 1332.75
$ws_CRP.Range("K122").Value = 5026.7142
$ws_CRP.Range("L122").Value = 3998.25
$ws_CRP.Range("M122").Value = -2576.7142
$ws_CRP.Range("N122").Value = -8898.25

# CUL row 34
$ws_CUL.Range("H34").Value = 1696.6666
$ws_CUL.Range("J34").Value = 2320
$ws_CUL.Range("L34").Value = 6960
$ws_CUL.Range("N34").Value = -7128

# CUL row 39
$ws_CUL.Range("H39").Value = 2899.5
$ws_CUL.Range("J39").Value = 2899.5
$ws_CUL.Range("L39").Value = 8698.5
$ws_CUL.Range("N39").Value = -9286.5

# CUL row 55
$ws_CUL.Range("H55").Value = 18835.666

# CUL row 75
$ws_CUL.Range("H75").Value = 1500
$ws_CUL.Range("I75").Value = 1000
$ws_CUL.Range("J75").Value = 1750
$ws_CUL.Range("K75").Value = 3000
$ws_CUL.Range("L75").Value = 5250
$ws_CUL.Range("M75").Value = -2002
$ws_CUL.Range("N75").Value = -7246

# CUL row 78
$ws_CUL.Range("H78").Value = 1500
$ws_CUL.Range("I78").Value = 1000
$ws_CUL.Range("J78").Value = 1750
$ws_CUL.Range("K78").Value = 9000
$ws_CUL.Range("L78").Value = 15750
$ws_CUL.Range("M78").Value = -4008
$ws_CUL.Range("N78").Value = -25734

# CUL row 92
$ws_CUL.Range("H92").Value = 300
$ws_CUL.Range("J92").Value = 300
$ws_CUL.Range("L92").Value = 900
$ws_CUL.Range("N92").Value = -3396

# CUL row 122
$ws_CUL.Range("H122").Value = 829.25
$ws_CUL.Range("I122").Value = 383
$ws_CUL.Range("K122").Value = 3447
$ws_CUL.Range("M122").Value = -997

# CUL row 131
$ws_CUL.Range("H131").Value = 5822853
$ws_CUL.Range("I131").Value = 125000536
$ws_CUL.Range("J131").Value = 9307.244000000001
$ws_CUL.Range("K131").Value = 375001608
$ws_CUL.Range("L131").Value = 27921.732
$ws_CUL.Range("M131").Value = -374996568
$ws_CUL.Range("N131").Value = -38001.732

# CUL row 137
$ws_CUL.Range("H137").Value = 5974.778
$ws_CUL.Range("I137").Value = 4560
$ws_CUL.Range("J137").Value = 7106.6
$ws_CUL.Range("K137").Value = 13680
$ws_CUL.Range("L137").Value = 21319.8
$ws_CUL.Range("M137").Value = -8580
$ws_CUL.Range("N137").Value = -31519.8

# GSM row 2
$ws_GSM.Range("H2").Value = 157.5
$ws_GSM.Range("I2").Value = 50
$ws_GSM.Range("K2").Value = 50
$ws_GSM.Range("M2").Value = 63

# GSM row 47
$ws_GSM.Range("H47").Value = 9000
$ws_GSM.Range("J47").Value = 9000
$ws_GSM.Range("L47").Value = 9000
$ws_GSM.Range("N47").Value = -10136

# LTW row 46
$ws_LTW.Range("H46").Value = 1636.3
$ws_LTW.Range("I46").Value = 1287.6666
$ws_LTW.Range("K46").Value = 1287.6666
$ws_LTW.Range("M46").Value = -1099.6666

# LTW row 93
$ws_LTW.Range("H93").Value = 592.25
$ws_LTW.Range("I93").Value = 208
$ws_LTW.Range("J93").Value = 1232.6666
$ws_LTW.Range("K93").Value = 208
$ws_LTW.Range("L93").Value = 1232.6666
$ws_LTW.Range("M93").Value = 1040
$ws_LTW.Range("N93").Value = -3728.6666

# LTW row 94
$ws_LTW.Range("H94").Value = 47790
$ws_LTW.Range("J94").Value = 47790
$ws_LTW.Range("L94").Value = 47790
$ws_LTW.Range("N94").Value = -49142

# WVR row 92
$ws_WVR.Range("H92").Value = 27782.5
$ws_WVR.Range("J92").Value = 27782.5
$ws_WVR.Range("L92").Value = 27782.5
$ws_WVR.Range("N92").Value = -32774.5

# WVR row 105
$ws_WVR.Range("H105").Value = 34997
$ws_WVR.Range("J105").Value = 34997
$ws_WVR.Range("L105").Value = 34997
$ws_WVR.Range("N105").Value = -41985

# WVR row 119
$ws_WVR.Range("H119").Value = 30000
$ws_WVR.Range("J119").Value = 30000
$ws_WVR.Range("L119").Value = 30000
$ws_WVR.Range("N119").Value = -39676

# WVR row 122
$ws_WVR.Range("H122").Value = 129631.1
$ws_WVR.Range("I122").Value = 183974
$ws_WVR.Range("J122").Value = 2831
$ws_WVR.Range("K122").Value = 551922
$ws_WVR.Range("L122").Value = 8493
$ws_WVR.Range("M122").Value = -549472
$ws_WVR.Range("N122").Value = -13393

# WVR row 132
$ws_WVR.Range("H132").Value = 2824.5518
$ws_WVR.Range("I132").Value = 2582.7083
$ws_WVR.Range("J132").Value = 3985.4
$ws_WVR.Range("K132").Value = 7748.124899999999
$ws_WVR.Range("L132").Value = 11956.2
$ws_WVR.Range("M132").Value = -5218.124899999999
$ws_WVR.Range("N132").Value = -17016.2

# WVR row 135
$ws_WVR.Range("H135").Value = 104178.75
$ws_WVR.Range("J135").Value = 104178.75
$ws_WVR.Range("L135").Value = 104178.75
$ws_WVR.Range("N135").Value = -114318.75
